$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# Plain numeric cells
$ws.Range("A$row").Value = 112544172
$ws.Range("B$row").Value = 90169

# Plain text cells
$ws.Range("C$row").Value = "Ovaliderad"
$ws.Range("D$row").Value = "LC"

# Numeric
$ws.Range("E$row").Value = 6031

# Plain text cells
$ws.Range("F$row").Value = "Blomkålssvamp"
$ws.Range("G$row").Value = "Sparassis crispa"
$ws.Range("H$row").Value = "(Wulfen:Fr.) Fr."

# "1" must be stored as TEXT, not a number -> use leading apostrophe
$ws.Range("I$row").Value = "'1"

$ws.Range("J$row").Value = "fruktkroppar"

# Empty-string text cells (still present as empty inline/shared strings)
$ws.Range("K$row").Value = "'"
$ws.Range("N$row").Value = "'"

$ws.Range("P$row").Value = "Lilla Bergsätter VNV 545 m, Ög"

# Numeric coordinates
$ws.Range("Q$row").Value = 563085
$ws.Range("R$row").Value = 6504264
$ws.Range("S$row").Value = 10

$ws.Range("T$row").Value = "Östergötland"
$ws.Range("U$row").Value = "Norrköping"
$ws.Range("V$row").Value = "Östergötland"
$ws.Range("W$row").Value = "Kvillinge"

# Dates stored as literal TEXT strings, not Excel date serials
$ws.Range("Y$row").Value = "'2023-10-05"
$ws.Range("AA$row").Value = "'2023-10-05"

# Boolean cells
$ws.Range("AD$row").Value = $false
$ws.Range("AE$row").Value = $false

# Empty-string text cell
$ws.Range("AF$row").Value = "'"

$ws.Range("AG$row").Value = $false

$ws.Range("AH$row").Value = "Barrskog"

# Empty-string text cell
$ws.Range("AT$row").Value = "'"

$ws.Range("AW$row").Value = "Mirjam Ideström"
$ws.Range("AX$row").Value = "Mirjam Ideström"

# Empty-string text cell
$ws.Range("AY$row").Value = "'"
